$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.024.73"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "'2.040.87"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "'247.24"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'56.28"
$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("E9").Value = "  -0.47%  "

$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "'15.93"
$ws.Range("E12").Value = "  -2.48%  "

$ws.Range("D13").Value = "'0.895"
$ws.Range("E13").Value = "  +11.96%  "

$ws.Range("D14").Value = "'2.337.69"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("D16").Value = "'2.044.54"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "'18.68"
$ws.Range("E17").Value = "  +12.13%  "

$ws.Range("D18").Value = "'37.078.00"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").Value = "'74.67"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'236.50"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  +4.29%  "

$ws.Range("D25").Value = "'171.32"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("E26").Value = "  +3.54%  "

$ws.Range("E27").Value = "  -8.29%  "

$ws.Range("D28").Value = "'20.07"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").Value = "'5.10"
$ws.Range("E30").Value = "  +8.62%  "

$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("D33").Value = "'4.64"
$ws.Range("E33").Value = "  +4.67%  "

$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").Value = "'0.0874"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  +5.71%  "

$ws.Range("E37").Value = "  +1.67%  "

$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("E39").Value = "  +9.63%  "

$ws.Range("E40").Value = "  +8.64%  "

$ws.Range("D41").Value = "'0.0991"
$ws.Range("E41").Value = "  -9.93%  "

$ws.Range("D42").Value = "'0.0223"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("D45").Value = "'97.60"

$ws.Range("E46").Value = "  -3.81%  "

$ws.Range("D47").Value = "'1.282.63"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("D50").Value = "'2.225.36"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "'44.40"
$ws.Range("E51").Value = "  +1.94%  "
